# Update "想去人数" (number of people interested) values for several events
# on the "展览" (Exhibitions) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 300
$ws1.Range("F4").Value = 8035
$ws1.Range("F5").Value = 5859
$ws1.Range("F7").Value = 87
$ws1.Range("F10").Value = 286
$ws1.Range("F11").Value = 385
$ws1.Range("F12").Value = 65

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 300
$ws4.Range("F4").Value = 8035
$ws4.Range("F5").Value = 5859
$ws4.Range("F7").Value = 87
$ws4.Range("F10").Value = 286
$ws4.Range("F14").Value = 385
$ws4.Range("F15").Value = 65
